$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "xx"
$ws.Range("B5").Value = "xvideo"
$ws.Range("C5").Value = "https://img-hw.xvideos-cdn.com/videos/thumbs169lll/1c/66/e1/1c66e1f0cfa69e0a380c294470e140ab/1c66e1f0cfa69e0a380c294470e140ab.4.jpg"

$ws.Range("D5").Select()
